$wb = $excel.ActiveWorkbook

# Loan RBI, Variable Instalments
# Insert a new (blank) column before column N on the "Repayment Schedule"
# sheet, shifting the old N:P columns (Late / Heading-Date / Outstanding)
# one column to the right, to O:Q.
$ws = $wb.Worksheets.Item("Repayment Schedule")
$ws.Columns("N").Insert()

# Make "Repayment Schedule" the active sheet/tab, with R7 selected.
$ws.Activate()
$ws.Range("R7").Select()
